# Apply an AutoFilter on the "State" column (C) of the Backlog sheet,
# filtering for "To do", then update C3's value to "In progress" without
# letting the filter re-hide/re-show rows (matching the commit's intent:
# "Stop the refreshing of value, to write something").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Turn on the AutoFilter for B2:C14 (Priority/State columns), filtering
# State == "To do".
$rng = $ws.Range("B2:C14")
$rng.AutoFilter(2, @("To do"))

# Excel records the filter range as a hidden workbook-scoped defined name
# on the sheet (_FilterDatabase).
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Backlog!`$B`$2:`$C`$14")
$fdb.Visible = $false

# Now change C3's state to "In progress" - this happens after the filter
# was already applied, so the row stays visible even though it no longer
# matches the filter criteria (matches the observed diff: row 3 keeps no
# hidden="1" attribute).
$ws.Range("C3").Value = "In progress"

# Move the active selection, as seen in the diff.
$ws.Range("C4").Select()
